$d = $word.ActiveDocument

$d.Content.Find.Execute("Play Jiggly Cash Slot for Free - Review 2021", $true, $false, $false, $false, $false, $true, 1, $false, "Play Jiggly Cash - Free Slot Game Review", 2)

$d.Content.Find.Execute("Cascading system", $true, $false, $false, $false, $false, $true, 1, $false, "Unique and visually stunning jelly candy theme", 2)
$d.Content.Find.Execute("Well-designed symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Cascading system creates exciting gameplay", 2)
$d.Content.Find.Execute("Exciting background and theme", $true, $false, $false, $false, $false, $true, 1, $false, "High winning potential and up to 10,000x bet", 2)
$d.Content.Find.Execute("10,000x potential win", $true, $false, $false, $false, $false, $true, 1, $false, "Accessible and compatible with all devices", 2)
$d.Content.Find.Execute("High volatility", $true, $false, $false, $false, $false, $true, 1, $false, "Bonus game can only be purchased at a cost", 2)
$d.Content.Find.Execute("Expensive bonus buy", $true, $false, $false, $false, $false, $true, 1, $false, "Volatility may not appeal to players seeking frequent wins", 2)

$d.Content.Find.Execute("Read our Jiggly Cash review and play this slot machine for free. Find out about gameplay features, visuals and theme, winning potential, and more.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Jiggly Cash, a free slot game with unique jelly candy theme and high winning potential.", 2)
